$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CR_REVIEW_001 ---
$ws.Range("A2").Value = "CR_REVIEW_001"
$ws.Range("B2").NumberFormat = "m/d"
$ws.Range("B2").Value = "17/10"
$ws.Range("C2").Value = "Jana Muhammed"
$ws.Range("D2").Value = "V1.2"
$ws.Range("E2").Value = "1//Document Title"
$ws.Range("F2").Value = "Document Title doesn't match the document's name"
$ws.Range("G2").Value = 'Document Title should be changed to ("PO_SAG_CR_Glasses")'
$ws.Range("H2").Value = "Closed"

# --- Row 3: CR_REVIEW_002 ---
$ws.Range("A3").Value = "CR_REVIEW_002"
$ws.Range("B3").NumberFormat = "m/d"
$ws.Range("B3").Value = "17/10"
$ws.Range("C3").Value = "Reham Essam"
$ws.Range("D3").Value = "V1.2"
$ws.Range("E3").Value = "1//Document Status/Author"
$ws.Range("F3").Value = "the author's name isn't the same as the last name in the document history table "
$ws.Range("G3").Value = "The author's name should be changed to the latest author's name "
$ws.Range("H3").Value = "Closed"

# --- Column widths for F:G ---
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 69.92

# --- Row 3 gains an explicit (custom) height flag, value unchanged ---
$ws.Rows.Item(3).RowHeight = 14.4

# --- Sheet view: scroll/selection state ---
$win = $excel.Application.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("H3").Select() | Out-Null

Write-Host "Edit complete"
